{"js": "// Update the date/title paragraph (first paragraph in the body, outside the table).\nconst title = context.document.body.paragraphs.getFirst();\ntitle.insertText(\"2026-02-13 Friday\", Word.InsertLocation.replace);\n\n// Update every cell of the math-drill table (20 rows x 5 columns) with the\n// new problems, preserving the existing table/cell formatting.\nconst table = context.document.body.tables.getFirst();\ntable.values = [\n  [\"46+35=\", \"80-6=\", \"27+15=\", \"45+29=\", \"71-44=\"],\n  [\"38+27=\", \"26+8=\", \"70-62=\", \"51-37=\", \"52-28=\"],\n  [\"16+38=\", \"14+68=\", \"59+24=\", \"11-2=\", \"61-59=\"],\n  [\"32+29=\", \"91-14=\", \"71-23=\", \"60-33=\", \"70-47=\"],\n  [\"90-51=\", \"68+19=\", \"27+14=\", \"25+36=\", \"2+79=\"],\n  [\"29+55=\", \"35+19=\", \"44+38=\", \"41-24=\", \"81-46=\"],\n  [\"13-9=\", \"35-18=\", \"9+85=\", \"56-19=\", \"91-5=\"],\n  [\"94-6=\", \"25+57=\", \"44+19=\", \"77+19=\", \"20-13=\"],\n  [\"96-9=\", \"22-3=\", \"18+46=\", \"48+34=\", \"59+13=\"],\n  [\"7+49=\", \"95-86=\", \"98-89=\", \"43-6=\", \"57+25=\"],\n  [\"17+14=\", \"50-12=\", \"39+6=\", \"91-9=\", \"59+28=\"],\n  [\"55+7=\", \"68-59=\", \"83-75=\", \"90-86=\", \"30-23=\"],\n  [\"26+38=\", \"57+38=\", \"47+19=\", \"58+35=\", \"40-11=\"],\n  [\"63-4=\", \"13+68=\", \"26+7=\", \"27+54=\", \"51-22=\"],\n  [\"29+6=\", \"90-9=\", \"98-19=\", \"6+87=\", \"12+49=\"],\n  [\"47+19=\", \"67-8=\", \"30-27=\", \"45+38=\", \"85-18=\"],\n  [\"67+8=\", \"8+74=\", \"28+5=\", \"43+9=\", \"46+5=\"],\n  [\"2+49=\", \"29+28=\", \"74-57=\", \"91-29=\", \"81-43=\"],\n  [\"45+17=\", \"39+35=\", \"48+46=\", \"5+89=\", \"71-24=\"],\n  [\"33-27=\", \"15+76=\", \"8+68=\", \"54-38=\", \"64-19=\"]\n];\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the title/date paragraph\n$d.Paragraphs.Item(1).Range.Text = '2026-02-13 Friday'\n\n# Update the table of math problems (20 rows x 5 columns), row-major order\n$tbl = $d.Tables.Item(1)\n$values = @(\n    '46+35=',\n    '80-6=',\n    '27+15=',\n    '45+29=',\n    '71-44=',\n    '38+27=',\n    '26+8=',\n    '70-62=',\n    '51-37=',\n    '52-28=',\n    '16+38=',\n    '14+68=',\n    '59+24=',\n    '11-2=',\n    '61-59=',\n    '32+29=',\n    '91-14=',\n    '71-23=',\n    '60-33=',\n    '70-47=',\n    '90-51=',\n    '68+19=',\n    '27+14=',\n    '25+36=',\n    '2+79=',\n    '29+55=',\n    '35+19=',\n    '44+38=',\n    '41-24=',\n    '81-46=',\n    '13-9=',\n    '35-18=',\n    '9+85=',\n    '56-19=',\n    '91-5=',\n    '94-6=',\n    '25+57=',\n    '44+19=',\n    '77+19=',\n    '20-13=',\n    '96-9=',\n    '22-3=',\n    '18+46=',\n    '48+34=',\n    '59+13=',\n    '7+49=',\n    '95-86=',\n    '98-89=',\n    '43-6=',\n    '57+25=',\n    '17+14=',\n    '50-12=',\n    '39+6=',\n    '91-9=',\n    '59+28=',\n    '55+7=',\n    '68-59=',\n    '83-75=',\n    '90-86=',\n    '30-23=',\n    '26+38=',\n    '57+38=',\n    '47+19=',\n    '58+35=',\n    '40-11=',\n    '63-4=',\n    '13+68=',\n    '26+7=',\n    '27+54=',\n    '51-22=',\n    '29+6=',\n    '90-9=',\n    '98-19=',\n    '6+87=',\n    '12+49=',\n    '47+19=',\n    '67-8=',\n    '30-27=',\n    '45+38=',\n    '85-18=',\n    '67+8=',\n    '8+74=',\n    '28+5=',\n    '43+9=',\n    '46+5=',\n    '2+49=',\n    '29+28=',\n    '74-57=',\n    '91-29=',\n    '81-43=',\n    '45+17=',\n    '39+35=',\n    '48+46=',\n    '5+89=',\n    '71-24=',\n    '33-27=',\n    '15+76=',\n    '8+68=',\n    '54-38=',\n    '64-19='\n)\n\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\nif (($rows * $cols) -ne $values.Count) {\n    throw \"Unexpected table size: $rows x $cols (expected $($values.Count) cells)\"\n}\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $values[$idx]\n        $idx++\n    }\n}\n\n"}
